$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 14.72626766666667
$ws.Cells.Item(2, 8).Value = 44.178803
$ws.Cells.Item(2, 9).Value = 0.08850173745156495
$ws.Cells.Item(2, 10).Value = 0.08850173745156498
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 28.06266266666667
$ws.Cells.Item(2, 14).Value = 84.187988
$ws.Cells.Item(2, 15).Value = 0.137091999334172
$ws.Cells.Item(2, 16).Value = 0.1370919993341719
$ws.Cells.Item(2, 17).Value = 413.2582818687071
$ws.Cells.Item(2, 18).Value = 3719.324536818364
$ws.Cells.Item(2, 19).Value = 0.012132880131783
$ws.Cells.Item(2, 20).Value = 0.01213288013178301

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 14.72626766666667
$ws.Cells.Item(3, 8).Value = 44.178803
$ws.Cells.Item(3, 9).Value = 0.08850173745156495
$ws.Cells.Item(3, 10).Value = 0.08850173745156498
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 41.27966333333334
$ws.Cells.Item(3, 14).Value = 123.83899
$ws.Cells.Item(3, 15).Value = 0.2016598227127667
$ws.Cells.Item(3, 16).Value = 0.2016598227127666
$ws.Cells.Item(3, 17).Value = 607.8953714365523
$ws.Cells.Item(3, 18).Value = 5471.058342928971
$ws.Cells.Item(3, 19).Value = 0.01784724468425441
$ws.Cells.Item(3, 20).Value = 0.01784724468425441

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 14.72626766666667
$ws.Cells.Item(4, 8).Value = 44.178803
$ws.Cells.Item(4, 9).Value = 0.08850173745156495
$ws.Cells.Item(4, 10).Value = 0.08850173745156498
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 127.2218856666667
$ws.Cells.Item(4, 14).Value = 381.665657
$ws.Cells.Item(4, 15).Value = 0.6215056237633366
$ws.Cells.Item(4, 16).Value = 0.6215056237633366
$ws.Cells.Item(4, 17).Value = 1873.503541385397
$ws.Cells.Item(4, 18).Value = 16861.53187246857
$ws.Cells.Item(4, 19).Value = 0.05500432753897393
$ws.Cells.Item(4, 20).Value = 0.05500432753897394

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 14.72626766666667
$ws.Cells.Item(5, 8).Value = 44.178803
$ws.Cells.Item(5, 9).Value = 0.08850173745156495
$ws.Cells.Item(5, 10).Value = 0.08850173745156498
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.135280666666667
$ws.Cells.Item(5, 14).Value = 24.405842
$ws.Cells.Item(5, 15).Value = 0.03974255418972485
$ws.Cells.Item(5, 16).Value = 0.03974255418972485
$ws.Cells.Item(5, 17).Value = 119.8023206407918
$ws.Cells.Item(5, 18).Value = 1078.220885767126
$ws.Cells.Item(5, 19).Value = 0.003517285096553622
$ws.Cells.Item(5, 20).Value = 0.003517285096553623

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 13.135456
$ws.Cells.Item(6, 8).Value = 39.406368
$ws.Cells.Item(6, 9).Value = 0.07894129758689367
$ws.Cells.Item(6, 10).Value = 0.07894129758689368
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 28.06266266666667
$ws.Cells.Item(6, 14).Value = 84.187988
$ws.Cells.Item(6, 15).Value = 0.137091999334172
$ws.Cells.Item(6, 16).Value = 0.1370919993341719
$ws.Cells.Item(6, 17).Value = 368.6158707008427
$ws.Cells.Item(6, 18).Value = 3317.542836307584
$ws.Cells.Item(6, 19).Value = 0.0108222203162211
$ws.Cells.Item(6, 20).Value = 0.0108222203162211

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 13.135456
$ws.Cells.Item(7, 8).Value = 39.406368
$ws.Cells.Item(7, 9).Value = 0.07894129758689367
$ws.Cells.Item(7, 10).Value = 0.07894129758689368
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 41.27966333333334
$ws.Cells.Item(7, 14).Value = 123.83899
$ws.Cells.Item(7, 15).Value = 0.2016598227127667
$ws.Cells.Item(7, 16).Value = 0.2016598227127666
$ws.Cells.Item(7, 17).Value = 542.2272014098133
$ws.Cells.Item(7, 18).Value = 4880.04481268832
$ws.Cells.Item(7, 19).Value = 0.01591928807608873
$ws.Cells.Item(7, 20).Value = 0.01591928807608873

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 13.135456
$ws.Cells.Item(8, 8).Value = 39.406368
$ws.Cells.Item(8, 9).Value = 0.07894129758689367
$ws.Cells.Item(8, 10).Value = 0.07894129758689368
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 127.2218856666667
$ws.Cells.Item(8, 14).Value = 381.665657
$ws.Cells.Item(8, 15).Value = 0.6215056237633366
$ws.Cells.Item(8, 16).Value = 0.6215056237633366
$ws.Cells.Item(8, 17).Value = 1671.117481411531
$ws.Cells.Item(8, 18).Value = 15040.05733270378
$ws.Cells.Item(8, 19).Value = 0.04906246039742953
$ws.Cells.Item(8, 20).Value = 0.04906246039742954

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 13.135456
$ws.Cells.Item(9, 8).Value = 39.406368
$ws.Cells.Item(9, 9).Value = 0.07894129758689367
$ws.Cells.Item(9, 10).Value = 0.07894129758689368
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 8.135280666666667
$ws.Cells.Item(9, 14).Value = 24.405842
$ws.Cells.Item(9, 15).Value = 0.03974255418972485
$ws.Cells.Item(9, 16).Value = 0.03974255418972485
$ws.Cells.Item(9, 17).Value = 106.8606212446507
$ws.Cells.Item(9, 18).Value = 961.745591201856
$ws.Cells.Item(9, 19).Value = 0.003137328797154317
$ws.Cells.Item(9, 20).Value = 0.003137328797154318

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 134.3676426666667
$ws.Cells.Item(10, 8).Value = 403.102928
$ws.Cells.Item(10, 9).Value = 0.807520962028172
$ws.Cells.Item(10, 10).Value = 0.8075209620281721
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 28.06266266666667
$ws.Cells.Item(10, 14).Value = 84.187988
$ws.Cells.Item(10, 15).Value = 0.137091999334172
$ws.Cells.Item(10, 16).Value = 0.1370919993341719
$ws.Cells.Item(10, 17).Value = 3770.713829469874
$ws.Cells.Item(10, 18).Value = 33936.42446522887
$ws.Cells.Item(10, 19).Value = 0.1107046631886961
$ws.Cells.Item(10, 20).Value = 0.1107046631886961

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 134.3676426666667
$ws.Cells.Item(11, 8).Value = 403.102928
$ws.Cells.Item(11, 9).Value = 0.807520962028172
$ws.Cells.Item(11, 10).Value = 0.8075209620281721
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 41.27966333333334
$ws.Cells.Item(11, 14).Value = 123.83899
$ws.Cells.Item(11, 15).Value = 0.2016598227127667
$ws.Cells.Item(11, 16).Value = 0.2016598227127666
$ws.Cells.Item(11, 17).Value = 5546.651052173637
$ws.Cells.Item(11, 18).Value = 49919.85946956273
$ws.Cells.Item(11, 19).Value = 0.162844534039444
$ws.Cells.Item(11, 20).Value = 0.162844534039444

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 134.3676426666667
$ws.Cells.Item(12, 8).Value = 403.102928
$ws.Cells.Item(12, 9).Value = 0.807520962028172
$ws.Cells.Item(12, 10).Value = 0.8075209620281721
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 127.2218856666667
$ws.Cells.Item(12, 14).Value = 381.665657
$ws.Cells.Item(12, 15).Value = 0.6215056237633366
$ws.Cells.Item(12, 16).Value = 0.6215056237633366
$ws.Cells.Item(12, 17).Value = 17094.50487263819
$ws.Cells.Item(12, 18).Value = 153850.5438537437
$ws.Cells.Item(12, 19).Value = 0.5018788192072887
$ws.Cells.Item(12, 20).Value = 0.5018788192072888

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 134.3676426666667
$ws.Cells.Item(13, 8).Value = 403.102928
$ws.Cells.Item(13, 9).Value = 0.807520962028172
$ws.Cells.Item(13, 10).Value = 0.8075209620281721
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 8.135280666666667
$ws.Cells.Item(13, 14).Value = 24.405842
$ws.Cells.Item(13, 15).Value = 0.03974255418972485
$ws.Cells.Item(13, 16).Value = 0.03974255418972485
$ws.Cells.Item(13, 17).Value = 1093.118485611709
$ws.Cells.Item(13, 18).Value = 9838.066370505376
$ws.Cells.Item(13, 19).Value = 0.03209294559274337
$ws.Cells.Item(13, 20).Value = 0.03209294559274338

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 4.165871666666667
$ws.Cells.Item(14, 8).Value = 12.497615
$ws.Cells.Item(14, 9).Value = 0.0250360029333692
$ws.Cells.Item(14, 10).Value = 0.02503600293336921
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 28.06266266666667
$ws.Cells.Item(14, 14).Value = 84.187988
$ws.Cells.Item(14, 15).Value = 0.137091999334172
$ws.Cells.Item(14, 16).Value = 0.1370919993341719
$ws.Cells.Item(14, 17).Value = 116.9054512942911
$ws.Cells.Item(14, 18).Value = 1052.14906164862
$ws.Cells.Item(14, 19).Value = 0.003432235697471777
$ws.Cells.Item(14, 20).Value = 0.003432235697471778

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 4.165871666666667
$ws.Cells.Item(15, 8).Value = 12.497615
$ws.Cells.Item(15, 9).Value = 0.0250360029333692
$ws.Cells.Item(15, 10).Value = 0.02503600293336921
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 41.27966333333334
$ws.Cells.Item(15, 14).Value = 123.83899
$ws.Cells.Item(15, 15).Value = 0.2016598227127667
$ws.Cells.Item(15, 16).Value = 0.2016598227127666
$ws.Cells.Item(15, 17).Value = 171.9657798898723
$ws.Cells.Item(15, 18).Value = 1547.69201900885
$ws.Cells.Item(15, 19).Value = 0.005048755912979539
$ws.Cells.Item(15, 20).Value = 0.005048755912979539

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 4.165871666666667
$ws.Cells.Item(16, 8).Value = 12.497615
$ws.Cells.Item(16, 9).Value = 0.0250360029333692
$ws.Cells.Item(16, 10).Value = 0.02503600293336921
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 127.2218856666667
$ws.Cells.Item(16, 14).Value = 381.665657
$ws.Cells.Item(16, 15).Value = 0.6215056237633366
$ws.Cells.Item(16, 16).Value = 0.6215056237633366
$ws.Cells.Item(16, 17).Value = 529.9900488786728
$ws.Cells.Item(16, 18).Value = 4769.910439908056
$ws.Cells.Item(16, 19).Value = 0.01556001661964435
$ws.Cells.Item(16, 20).Value = 0.01556001661964435

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 4.165871666666667
$ws.Cells.Item(17, 8).Value = 12.497615
$ws.Cells.Item(17, 9).Value = 0.0250360029333692
$ws.Cells.Item(17, 10).Value = 0.02503600293336921
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 8.135280666666667
$ws.Cells.Item(17, 14).Value = 24.405842
$ws.Cells.Item(17, 15).Value = 0.03974255418972485
$ws.Cells.Item(17, 16).Value = 0.03974255418972485
$ws.Cells.Item(17, 17).Value = 33.89053522964778
$ws.Cells.Item(17, 18).Value = 305.01481706683
$ws.Cells.Item(17, 19).Value = 0.0009949947032735359
$ws.Cells.Item(17, 20).Value = 0.0009949947032735361
